# Update TPM-derived NATMI ligand-receptor metrics (Sema6d-Kdr) with newly computed values.
# Only numeric value cells change; no rows/columns are inserted or removed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 58.62893533333334
$ws.Range("H2").Value = 175.886806
$ws.Range("I2").Value = 0.5702456571409142
$ws.Range("J2").Value = 0.5702456571409142
$ws.Range("M2").Value = 184.1023456666667
$ws.Range("N2").Value = 552.307037
$ws.Range("O2").Value = 0.9813423747591566
$ws.Range("P2").Value = 0.9813423747591565
$ws.Range("Q2").Value = 10793.72451880598
$ws.Range("R2").Value = 97143.52066925383
$ws.Range("S2").Value = 0.5596062273747606
$ws.Range("T2").Value = 0.5596062273747605

# Row 3
$ws.Range("G3").Value = 58.62893533333334
$ws.Range("H3").Value = 175.886806
$ws.Range("I3").Value = 0.5702456571409142
$ws.Range("J3").Value = 0.5702456571409142
$ws.Range("O3").Value = 0.002303378255889225
$ws.Range("P3").Value = 0.002303378255889224
$ws.Range("Q3").Value = 25.33471599326156
$ws.Range("R3").Value = 228.012443939354
$ws.Range("S3").Value = 0.001313491447173644
$ws.Range("T3").Value = 0.001313491447173644

# Row 4
$ws.Range("G4").Value = 58.62893533333334
$ws.Range("H4").Value = 175.886806
$ws.Range("I4").Value = 0.5702456571409142
$ws.Range("J4").Value = 0.5702456571409142
$ws.Range("M4").Value = 1.367901
$ws.Range("N4").Value = 4.103703
$ws.Range("O4").Value = 0.007291483500193526
$ws.Range("P4").Value = 0.007291483500193526
$ws.Range("Q4").Value = 80.19857927140201
$ws.Range("R4").Value = 721.787213442618
$ws.Range("S4").Value = 0.004157936800099991
$ws.Range("T4").Value = 0.004157936800099991

# Row 5
$ws.Range("G5").Value = 58.62893533333334
$ws.Range("H5").Value = 175.886806
$ws.Range("I5").Value = 0.5702456571409142
$ws.Range("J5").Value = 0.5702456571409142
$ws.Range("M5").Value = 1.700197666666667
$ws.Range("N5").Value = 5.100593
$ws.Range("O5").Value = 0.009062763484760617
$ws.Range("P5").Value = 0.009062763484760615
$ws.Range("Q5").Value = 99.68077905288423
$ws.Range("R5").Value = 897.127011475958
$ws.Range("S5").Value = 0.00516800151888
$ws.Range("T5").Value = 0.005168001518879999

# Row 6
$ws.Range("I6").Value = 0.1389799721218762
$ws.Range("J6").Value = 0.1389799721218763
$ws.Range("M6").Value = 184.1023456666667
$ws.Range("N6").Value = 552.307037
$ws.Range("O6").Value = 0.9813423747591566
$ws.Range("P6").Value = 0.9813423747591565
$ws.Range("Q6").Value = 2630.640871928942
$ws.Range("R6").Value = 23675.76784736048
$ws.Range("S6").Value = 0.1363869358860434
$ws.Range("T6").Value = 0.1363869358860434

# Row 7
$ws.Range("I7").Value = 0.1389799721218762
$ws.Range("J7").Value = 0.1389799721218763
$ws.Range("O7").Value = 0.002303378255889225
$ws.Range("P7").Value = 0.002303378255889224
$ws.Range("S7").Value = 0.0003201234457896203
$ws.Range("T7").Value = 0.0003201234457896204

# Row 8
$ws.Range("I8").Value = 0.1389799721218762
$ws.Range("J8").Value = 0.1389799721218763
$ws.Range("M8").Value = 1.367901
$ws.Range("N8").Value = 4.103703
$ws.Range("O8").Value = 0.007291483500193526
$ws.Range("P8").Value = 0.007291483500193526
$ws.Range("Q8").Value = 19.54595562768
$ws.Range("R8").Value = 175.91360064912
$ws.Range("S8").Value = 0.001013370173584017
$ws.Range("T8").Value = 0.001013370173584017

# Row 9
$ws.Range("I9").Value = 0.1389799721218762
$ws.Range("J9").Value = 0.1389799721218763
$ws.Range("M9").Value = 1.700197666666667
$ws.Range("N9").Value = 5.100593
$ws.Range("O9").Value = 0.009062763484760617
$ws.Range("P9").Value = 0.009062763484760615
$ws.Range("Q9").Value = 24.29414712830222
$ws.Range("R9").Value = 218.64732415472
$ws.Range("S9").Value = 0.001259542616459188
$ws.Range("T9").Value = 0.001259542616459188

# Row 10
$ws.Range("G10").Value = 27.27518533333334
$ws.Range("H10").Value = 81.825556
$ws.Range("I10").Value = 0.265288050953297
$ws.Range("J10").Value = 0.2652880509532971
$ws.Range("M10").Value = 184.1023456666667
$ws.Range("N10").Value = 552.307037
$ws.Range("O10").Value = 0.9813423747591566
$ws.Range("P10").Value = 0.9813423747591565
$ws.Range("Q10").Value = 5021.425598359731
$ws.Range("R10").Value = 45192.83038523758
$ws.Range("S10").Value = 0.2603384059177367
$ws.Range("T10").Value = 0.2603384059177367

# Row 11
$ws.Range("G11").Value = 27.27518533333334
$ws.Range("H11").Value = 81.825556
$ws.Range("I11").Value = 0.265288050953297
$ws.Range("J11").Value = 0.2652880509532971
$ws.Range("O11").Value = 0.002303378255889225
$ws.Range("P11").Value = 0.002303378255889224
$ws.Range("Q11").Value = 11.78614399451156
$ws.Range("R11").Value = 106.075295950604
$ws.Range("S11").Value = 0.0006110587281130572
$ws.Range("T11").Value = 0.0006110587281130572

# Row 12
$ws.Range("G12").Value = 27.27518533333334
$ws.Range("H12").Value = 81.825556
$ws.Range("I12").Value = 0.265288050953297
$ws.Range("J12").Value = 0.2652880509532971
$ws.Range("M12").Value = 1.367901
$ws.Range("N12").Value = 4.103703
$ws.Range("O12").Value = 0.007291483500193526
$ws.Range("P12").Value = 0.007291483500193526
$ws.Range("Q12").Value = 37.309753292652
$ws.Range("R12").Value = 335.7877796338681
$ws.Range("S12").Value = 0.001934343446324465
$ws.Range("T12").Value = 0.001934343446324465

# Row 13
$ws.Range("G13").Value = 27.27518533333334
$ws.Range("H13").Value = 81.825556
$ws.Range("I13").Value = 0.265288050953297
$ws.Range("J13").Value = 0.2652880509532971
$ws.Range("M13").Value = 1.700197666666667
$ws.Range("N13").Value = 5.100593
$ws.Range("O13").Value = 0.009062763484760617
$ws.Range("P13").Value = 0.009062763484760615
$ws.Range("Q13").Value = 46.37320646163423
$ws.Range("R13").Value = 417.358858154708
$ws.Range("S13").Value = 0.002404242861122854
$ws.Range("T13").Value = 0.002404242861122854

# Row 14
$ws.Range("E14").Value = 3.0
$ws.Range("F14").Value = 1.0
$ws.Range("G14").Value = 2.620337
$ws.Range("H14").Value = 7.861011
$ws.Range("I14").Value = 0.02548631978391236
$ws.Range("J14").Value = 0.02548631978391236
$ws.Range("M14").Value = 184.1023456666667
$ws.Range("N14").Value = 552.307037
$ws.Range("O14").Value = 0.9813423747591566
$ws.Range("P14").Value = 0.9813423747591565
$ws.Range("Q14").Value = 482.4101881371563
$ws.Range("R14").Value = 4341.691693234407
$ws.Range("S14").Value = 0.02501080558061583
$ws.Range("T14").Value = 0.02501080558061583

# Row 15
$ws.Range("E15").Value = 3.0
$ws.Range("F15").Value = 1.0
$ws.Range("G15").Value = 2.620337
$ws.Range("H15").Value = 7.861011
$ws.Range("I15").Value = 0.02548631978391236
$ws.Range("J15").Value = 0.02548631978391236
$ws.Range("O15").Value = 0.002303378255889225
$ws.Range("P15").Value = 0.002303378255889224
$ws.Range("Q15").Value = 1.132299150994333
$ws.Range("R15").Value = 10.190692358949
$ws.Range("S15").Value = 0.00005870463481290308
$ws.Range("T15").Value = 0.00005870463481290308

# Row 16
$ws.Range("E16").Value = 3.0
$ws.Range("F16").Value = 1.0
$ws.Range("G16").Value = 2.620337
$ws.Range("H16").Value = 7.861011
$ws.Range("I16").Value = 0.02548631978391236
$ws.Range("J16").Value = 0.02548631978391236
$ws.Range("M16").Value = 1.367901
$ws.Range("N16").Value = 4.103703
$ws.Range("O16").Value = 0.007291483500193526
$ws.Range("P16").Value = 0.007291483500193526
$ws.Range("Q16").Value = 3.584361602637
$ws.Range("R16").Value = 32.259254423733
$ws.Range("S16").Value = 0.0001858330801850528
$ws.Range("T16").Value = 0.0001858330801850528

# Row 17
$ws.Range("E17").Value = 3.0
$ws.Range("F17").Value = 1.0
$ws.Range("G17").Value = 2.620337
$ws.Range("H17").Value = 7.861011
$ws.Range("I17").Value = 0.02548631978391236
$ws.Range("J17").Value = 0.02548631978391236
$ws.Range("M17").Value = 1.700197666666667
$ws.Range("N17").Value = 5.100593
$ws.Range("O17").Value = 0.009062763484760617
$ws.Range("P17").Value = 0.009062763484760615
$ws.Range("Q17").Value = 4.455090853280333
$ws.Range("R17").Value = 40.09581767952299
$ws.Range("S17").Value = 0.000230976488298573
$ws.Range("T17").Value = 0.000230976488298573
